$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 25.06781033333333
$ws.Range("H2").Value = 75.203431
$ws.Range("I2").Value = 0.7308832858982242
$ws.Range("J2").Value = 0.7308832858982242
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.10446833333334
$ws.Range("N2").Value = 102.313405
$ws.Range("O2").Value = 0.5118942073015388
$ws.Range("P2").Value = 0.5118942073015389
$ws.Range("Q2").Value = 854.9243436991727
$ws.Range("R2").Value = 7694.319093292554
$ws.Range("S2").Value = 0.3741349202648154
$ws.Range("T2").Value = 0.3741349202648155

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 25.06781033333333
$ws.Range("H3").Value = 75.203431
$ws.Range("I3").Value = 0.7308832858982242
$ws.Range("J3").Value = 0.7308832858982242
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 31.083557
$ws.Range("N3").Value = 93.25067100000001
$ws.Range("O3").Value = 0.4665515560925921
$ws.Range("P3").Value = 0.4665515560925922
$ws.Range("Q3").Value = 779.1967113613556
$ws.Range("R3").Value = 7012.770402252201
$ws.Range("S3").Value = 0.3409947343578834
$ws.Range("T3").Value = 0.3409947343578835

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 25.06781033333333
$ws.Range("H4").Value = 75.203431
$ws.Range("I4").Value = 0.7308832858982242
$ws.Range("J4").Value = 0.7308832858982242
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.436030666666667
$ws.Range("N4").Value = 4.308092
$ws.Range("O4").Value = 0.02155423660586901
$ws.Range("P4").Value = 0.02155423660586901
$ws.Range("Q4").Value = 35.99814438485022
$ws.Range("R4").Value = 323.983299463652
$ws.Range("S4").Value = 0.01575363127552533
$ws.Range("T4").Value = 0.01575363127552533

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.8930513333333332
$ws.Range("H5").Value = 2.679154
$ws.Range("I5").Value = 0.02603802583086097
$ws.Range("J5").Value = 0.02603802583086097
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.10446833333334
$ws.Range("N5").Value = 102.313405
$ws.Range("O5").Value = 0.5118942073015388
$ws.Range("P5").Value = 0.5118942073015389
$ws.Range("Q5").Value = 30.45704091770778
$ws.Range("R5").Value = 274.11336825937
$ws.Range("S5").Value = 0.01332871459238557
$ws.Range("T5").Value = 0.01332871459238557

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.8930513333333332
$ws.Range("H6").Value = 2.679154
$ws.Range("I6").Value = 0.02603802583086097
$ws.Range("J6").Value = 0.02603802583086097
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.083557
$ws.Range("N6").Value = 93.25067100000001
$ws.Range("O6").Value = 0.4665515560925921
$ws.Range("P6").Value = 0.4665515560925922
$ws.Range("Q6").Value = 27.75921202359266
$ws.Range("R6").Value = 249.832908212334
$ws.Range("S6").Value = 0.0121480814689673
$ws.Range("T6").Value = 0.0121480814689673

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.8930513333333332
$ws.Range("H7").Value = 2.679154
$ws.Range("I7").Value = 0.02603802583086097
$ws.Range("J7").Value = 0.02603802583086097
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.436030666666667
$ws.Range("N7").Value = 4.308092
$ws.Range("O7").Value = 0.02155423660586901
$ws.Range("P7").Value = 0.02155423660586901
$ws.Range("Q7").Value = 1.282449101574222
$ws.Range("R7").Value = 11.542041914168
$ws.Range("S7").Value = 0.0005612297695081063
$ws.Range("T7").Value = 0.0005612297695081063

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.337104666666667
$ws.Range("H8").Value = 25.011314
$ws.Range("I8").Value = 0.2430786882709149
$ws.Range("J8").Value = 0.2430786882709149
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.10446833333334
$ws.Range("N8").Value = 102.313405
$ws.Range("O8").Value = 0.5118942073015388
$ws.Range("P8").Value = 0.5118942073015389
$ws.Range("Q8").Value = 284.3325220960189
$ws.Range("R8").Value = 2558.99269886417
$ws.Range("S8").Value = 0.1244305724443378
$ws.Range("T8").Value = 0.1244305724443379

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.337104666666667
$ws.Range("H9").Value = 25.011314
$ws.Range("I9").Value = 0.2430786882709149
$ws.Range("J9").Value = 0.2430786882709149
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 31.083557
$ws.Range("N9").Value = 93.25067100000001
$ws.Range("O9").Value = 0.4665515560925921
$ws.Range("P9").Value = 0.4665515560925922
$ws.Range("Q9").Value = 259.1468681212993
$ws.Range("R9").Value = 2332.321813091694
$ws.Range("S9").Value = 0.1134087402657415
$ws.Range("T9").Value = 0.1134087402657415

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.337104666666667
$ws.Range("H10").Value = 25.011314
$ws.Range("I10").Value = 0.2430786882709149
$ws.Range("J10").Value = 0.2430786882709149
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.436030666666667
$ws.Range("N10").Value = 4.308092
$ws.Range("O10").Value = 0.02155423660586901
$ws.Range("P10").Value = 0.02155423660586901
$ws.Range("Q10").Value = 11.97233797254311
$ws.Range("R10").Value = 107.751041752888
$ws.Range("S10").Value = 0.005239375560835575
$ws.Range("T10").Value = 0.005239375560835575
